# "DB uniqueness checker added": a duplicate visit record for the
# 100160016W / Abcfinance*20 account was collapsed, so its visit count
# drops by one and its "last visit" timestamp moves to the next
# (now-latest) recorded visit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4
$ws.Range("D2").Value = "2020-03-07 00:13:44.591811"
